# 9.3.1 worksheet update: add a new "2021" data column (O), correct a couple
# of previously entered data points, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new year column (O), cloning the formatting already used for
#     the existing "2020" column (N) for the header-border row (3),
#     the year-label row (4) and the data row (5). ---
$ws.Range("N3:N5").Copy()
$ws.Range("O3:O5").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $False

# New year label for 2021
$ws.Range("O4").Value = 2021

# Corrected historical value for 2018 (L5) and 2020 (N5), plus the brand new
# 2021 data point (O5)
$ws.Range("L5").Value = 1.6
$ws.Range("N5").Value = 3.1
$ws.Range("O5").Value = 4.1

# Move the active cell/selection as recorded in the saved view
$ws.Range("P4").Select()
